$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A15").Value = 7
$ws.Range("B15").Value = 542
$ws.Range("A16").Value = 15
$ws.Range("B16").Value = 331
$ws.Range("A17").Value = 21
$ws.Range("B17").Value = 240
